# Remove os índices da planilha a ser exportada
# The first column (A) of the sheet contains a numeric row index (0,1,2,...)
# which is no longer wanted in the exported spreadsheet. Delete that whole
# column so every other column shifts one position to the left
# (B->A, C->B, D->C, E->D, F->E, G->F, H->G), and refresh the quotation
# ("Cotação") values together with the figures that derive from them
# ("Preço de Compra" and "Preço de Venda").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the index column (old column A) - everything shifts left.
$ws.Columns.Item(1).Delete()

# Updated currency quotations (column D after the shift) and the values
# that depend on them (columns E and G), recalculated with the new rates.
# Row 2 - Câmera Canon (Dólar)
$ws.Range("D2").Value = 5.233274
$ws.Range("E2").Value = 5233.22166726
$ws.Range("G2").Value = 7326.510334163999

# Row 3 - Carro Renault (Euro)
$ws.Range("D3").Value = 5.5827
$ws.Range("E3").Value = 25122.15
$ws.Range("G3").Value = 50244.3

# Row 4 - Notebook Dell (Dólar)
$ws.Range("D4").Value = 5.233274
$ws.Range("E4").Value = 4709.89426726
$ws.Range("G4").Value = 8006.820254341999

# Row 5 - IPhone (Dólar)
$ws.Range("D5").Value = 5.233274
$ws.Range("E5").Value = 4181.385926
$ws.Range("G5").Value = 7108.3560742

# Row 6 - Carro Fiat (Euro)
$ws.Range("D6").Value = 5.5827
$ws.Range("E6").Value = 16748.1
$ws.Range("G6").Value = 31821.39

# Row 7 - Celular Xiaomi (Dólar)
$ws.Range("D7").Value = 5.233274
$ws.Range("E7").Value = 2514.48349152
$ws.Range("G7").Value = 5028.96698304

# Row 8 - Joia 20g (Ouro)
$ws.Range("D8").Value = 310.42
$ws.Range("E8").Value = 6208.400000000001
$ws.Range("G8").Value = 7139.66
